$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$featureNames = @(
    'deferi',
    'aumento',
    'deferimento liminar',
    'arquivado definitivo',
    'assessoria',
    'decisão proferida',
    'opinou concessão',
    'medidas cautelares',
    'concessão ordem',
    'deferida',
    'senha relatório',
    'decretada desfavor',
    'código senha',
    'justiça nº',
    'contornos',
    'arquivado',
    'formalizada superior',
    'cautelares previstas',
    'julgado turma',
    'opina deferimento',
    'ordem prisão',
    'aurélio assessor',
    'relatório http',
    'liminar espécie',
    'liberdade restritiva',
    'república concessão',
    'informado',
    'liminar assessor',
    'liminar deferida',
    'violência grave',
    'relativização',
    'liminar hc',
    'liminar suspender',
    'manifesta',
    'república parecer',
    'ministro gilson',
    'mulheres',
    'referida decisão',
    'opina concessão',
    'outro motivo',
    'paciente decretada',
    'proferida ministro',
    'indispensável custódia',
    'segundo autos',
    'implicou deferimento',
    'sob argumento',
    'assessoria prestou',
    'sumária',
    'substituir',
    'stj indeferiu',
    'causas aumento',
    'circunstâncias favoráveis',
    'concessão',
    'sobrestamento',
    'cpp art',
    'sob código',
    'deferi pedido',
    'deferimento',
    'deserção',
    'efeitos ordem',
    'ser julgado',
    'senha primeira',
    'previstas art',
    'restritiva direitos',
    'procuradoria geral',
    'restritiva',
    'senha',
    'revelou contornos',
    'resumida prisão',
    'resumida',
    'suspender',
    'suspender efeitos',
    'análise pedido',
    'preventiva fundamentos',
    'ficou',
    'assim resumida',
    'assim revelou',
    'aurélio decisão',
    'brasília residência',
    'campo precário',
    'contornos impetração',
    'decisão implicou',
    'deferida assessoria',
    'efêmero',
    'eis informado',
    'espécie ficou',
    'ficou assim',
    'precário efêmero',
    'fundamentos insubsistência',
    'gabinete prestou',
    'http',
    'http sob',
    'impetração eis',
    'indeferiu liminarmente',
    'informado análise',
    'informações paciente',
    'manifesta ilegalidade',
    'nº ministro',
    'precário',
    'jurisprudenciais'
)

$importanceValues = @(
    0.05,
    0.05,
    0.03,
    0.03,
    0.03,
    0.03,
    0.03,
    0.03,
    0.02,
    0.02,
    0.02,
    0.02,
    0.02,
    0.02,
    0.02,
    0.02,
    0.02,
    0.02,
    0.02,
    0.02,
    0.02,
    0.02,
    0.02,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0.01,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0
)

for ($i = 0; $i -lt $featureNames.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $featureNames[$i]
    $ws.Cells.Item($row, 2).Value = $importanceValues[$i]
}

